$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (column headers) shrinks a bit
$ws.Rows(2).RowHeight = 75

# Matrix body (B3:P9): center-align the bordered, empty tracking cells
$ws.Range("B3:P9").HorizontalAlignment = -4108

# Fill in the traceability matrix "X" marks
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("G3").Value = "X"
$ws.Range("I5").Value = "X"
$ws.Range("J6").Value = "X"
$ws.Range("L7").Value = "X"

# Selection cursor moved
$ws.Range("J17").Select()
